$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A256").Value2 = 44330
$ws.Range("B256").Value2 = 0
$ws.Range("C256").Value2 = 5
$ws.Range("D256").Value2 = 152.6717557251908

$ws.Range("A257").Value2 = 44331
$ws.Range("B257").Value2 = 1
$ws.Range("C257").Value2 = 3
$ws.Range("D257").Value2 = 91.6030534351145

$ws.Range("A258").Value2 = 44332
$ws.Range("B258").Value2 = 1
$ws.Range("C258").Value2 = 3
$ws.Range("D258").Value2 = 91.6030534351145

$ws.Range("A259").Value2 = 44333
$ws.Range("B259").Value2 = 0
$ws.Range("C259").Value2 = 3
$ws.Range("D259").Value2 = 91.6030534351145

$ws.Range("A260").Value2 = 44334
$ws.Range("B260").Value2 = 0
$ws.Range("C260").Value2 = 3
$ws.Range("D260").Value2 = 91.6030534351145

$ws.Range("A261").Value2 = 44335
$ws.Range("B261").Value2 = 0
$ws.Range("C261").Value2 = 3
$ws.Range("D261").Value2 = 91.6030534351145

$ws.Range("A262").Value2 = 44336
$ws.Range("B262").Value2 = 0
$ws.Range("C262").Value2 = 2
$ws.Range("D262").Value2 = 61.06870229007634

$ws.Range("A263").Value2 = 44337
$ws.Range("B263").Value2 = 0
$ws.Range("C263").Value2 = 2
$ws.Range("D263").Value2 = 61.06870229007634

$ws.Range("A264").Value2 = 44338
$ws.Range("B264").Value2 = 0
$ws.Range("C264").Value2 = 1
$ws.Range("D264").Value2 = 30.53435114503817

$ws.Range("A265").Value2 = 44339
$ws.Range("B265").Value2 = 0
$ws.Range("C265").Value2 = 0
$ws.Range("D265").Value2 = 0

$ws.Range("A266").Value2 = 44340
$ws.Range("B266").Value2 = 0
$ws.Range("C266").Value2 = 0
$ws.Range("D266").Value2 = 0

$ws.Range("A267").Value2 = 44341
$ws.Range("B267").Value2 = 0
$ws.Range("C267").Value2 = 0
$ws.Range("D267").Value2 = 0

$ws.Range("A268").Value2 = 44342
$ws.Range("B268").Value2 = 0
$ws.Range("C268").Value2 = 0
$ws.Range("D268").Value2 = 0

$ws.Range("A269").Value2 = 44343
$ws.Range("B269").Value2 = 0
$ws.Range("C269").Value2 = 0
$ws.Range("D269").Value2 = 0

$ws.Range("A255").Copy()
$ws.Range("A256:A269").PasteSpecial(-4122)
$excel.CutCopyMode = 0

